$d = $word.ActiveDocument

# --- Remove the "_GoBack" bookmark that currently sits in front of the
#     "Power Failure Management" run (paragraph 1). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Insert a brand-new paragraph right after paragraph 1. It inherits the
#     centered / 32pt formatting from paragraph 1's paragraph mark. ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter() | Out-Null

$p2 = $d.Paragraphs.Item(2)

# Type the new paragraph's text. A trailing placeholder character is added
# so that the bookmark anchor point below does not land exactly on the
# paragraph mark (doing so makes the host collapse/snap the bookmark to the
# whole first paragraph instead of the intended position); the placeholder
# is stripped again immediately afterwards.
$p2.Range.Text = "ljmknkknjbjbjX"

$anchorPos = $p2.Range.End - 2   # position right after "...jbjbj", before the "X"
$bmRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character again; the bookmark (anchored just
# before it) stays put, ending up right after the real text and before the
# paragraph mark - exactly where it needs to be.
$placeholderRange = $d.Range($anchorPos, $anchorPos + 1)
$placeholderRange.Delete()
